$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.841744661331177
$ws.Range("B1").Value = 6.087389469146729
$ws.Range("C1").Value = 2.24291467666626
$ws.Range("D1").Value = 1.469248294830322
$ws.Range("E1").Value = 1.208699107170105
